# Actualizado a 24 de Marzo de 2020
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 22) with the case counts for 24 March 2020.
$newRow = 22

# Copy the formatting from the row above (row 21) down into the new row,
# so the new row inherits the same number formats / styles (date format
# in column A, general number format elsewhere).
$ws.Range("A21:S21").Copy()
$ws.Range("A$newRow`:S$newRow").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the values for the new row.
$ws.Cells.Item($newRow, 1).Value = 43913  # fecha (24/03/2020)
$ws.Cells.Item($newRow, 2).Value = 21     # dia
$ws.Cells.Item($newRow, 3).Value = 2      # Arica y Parinacota
$ws.Cells.Item($newRow, 4).Value = 4      # Tarapaca
$ws.Cells.Item($newRow, 5).Value = 13     # Antofagasta
$ws.Cells.Item($newRow, 6).Value = 1      # Atacama
$ws.Cells.Item($newRow, 7).Value = 11     # Coquimbo
$ws.Cells.Item($newRow, 8).Value = 25     # Valparaiso
$ws.Cells.Item($newRow, 9).Value = 540    # Metropolitana
$ws.Cells.Item($newRow, 10).Value = 9     # O'Higgins
$ws.Cells.Item($newRow, 11).Value = 29    # Maule
$ws.Cells.Item($newRow, 12).Value = 105   # Nuble
$ws.Cells.Item($newRow, 13).Value = 73    # Biobio
$ws.Cells.Item($newRow, 14).Value = 59    # Araucania
$ws.Cells.Item($newRow, 15).Value = 6     # Los Rios
$ws.Cells.Item($newRow, 16).Value = 36    # Los Lagos
$ws.Cells.Item($newRow, 17).Value = 1     # Aysen
$ws.Cells.Item($newRow, 18).Value = 8     # Magallanes
$ws.Cells.Item($newRow, 19).Value = 922   # total

# Update the selection to reflect where the user left off editing.
$ws.Range("T22").Select()
